# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# This updates the G column (header "K") values for rows 2-28 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 4
    13 = 5
    14 = 5
    15 = 7
    16 = 7
    17 = 3
    18 = 8
    19 = 5
    20 = 2
    21 = 7
    22 = 3
    23 = 3
    24 = 5
    25 = 7
    26 = 4
    27 = 5
    28 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
